$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 26: "Create Contact" API entry ---------------------------------

$jsonData = "{`n    ""user_id"" : 1,`n    ""agent_id"" : 10,`n    ""first_name"" : ""Junaid"",`n    ""email"" : ""junaid.ansari@1point1.com"",`n    ""phone_number"" : ""9821209237""`n}"

$curlData = "curl --location 'http://1msg.1point1.in:3001/api/chat/bot/create/contact/' \`n--header 'Content-Type: application/json' \`n--data-raw '{`n    ""user_id"" : 1,`n    ""agent_id"" : 10,`n    ""first_name"" : ""Junaid"",`n    ""email"" : ""junaid.ansari@1point1.com"",`n    ""phone_number"" : ""9821209237""`n}'"

$ws.Range("A26").Value = 12
$ws.Range("B26").Value = "Create Contact"
$ws.Range("C26").Value = "POST"
$ws.Range("D26").Value = "http://1msg.1point1.in:3001/api/chat/bot/create/contact/"
$ws.Range("F26").Value = $jsonData
$ws.Range("G26").Value = $curlData
$ws.Range("H26").Value = "Create Contact in Chat"

# Hyperlink on the endpoint URL, matching the style used by the other rows
$ws.Hyperlinks.Add($ws.Range("D26"), "http://1msg.1point1.in:3001/api/chat/bot/create/contact/")
$ws.Range("D3").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height to accommodate the wrapped content
$ws.Rows.Item(26).RowHeight = 230.4

# --- View state ---------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F26").Select()
